$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B63").Value = 237
$ws.Range("B64").Value = 237
$ws.Range("B65").Value = 237

$ws.Range("B69").Value = 298
$ws.Range("B70").Value = 298
$ws.Range("B71").Value = 298

$ws.Range("B72").Value = 233
$ws.Range("B73").Value = 233
$ws.Range("B74").Value = 233

$ws.Range("B75").Value = 219
$ws.Range("B76").Value = 219
$ws.Range("B77").Value = 219

$ws.Range("B78").Value = 229
$ws.Range("B79").Value = 229
$ws.Range("B80").Value = 229

$ws.Range("B81").Value = 216
$ws.Range("B82").Value = 216
$ws.Range("B83").Value = 216

$ws.Range("B84").Value = 271
$ws.Range("B85").Value = 271
$ws.Range("B86").Value = 271

$ws.Range("B87").Value = 215
$ws.Range("B88").Value = 215
$ws.Range("B89").Value = 215

$ws.Range("B90").Value = 298
$ws.Range("B91").Value = 298
$ws.Range("B92").Value = 298

$ws.Range("B93").Value = 302
$ws.Range("B94").Value = 302
$ws.Range("B95").Value = 302
